# Apply crypto price/volume updates to the active worksheet.
# Generated from the authoritative diff of the Price (D) and Volume(1h) (E)
# columns for rows 2-51. All values in these columns are stored as plain
# text in the workbook (e.g. "0.5325", "  +1.32%  "), matching the source
# data feed's formatting. Cells whose new text would otherwise be
# auto-recognized by Excel as a number (e.g. "301.13") have their number
# format forced to Text ("@") first so the literal string is preserved,
# exactly like the existing text-valued cells in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textUpdates = @{
    'D5' = '301.13'
    'D7' = '0.5325'
    'D8' = '0.3742'
    'D9' = '0.07181'
    'D10' = '21.61'
    'D11' = '0.8886'
    'D12' = '0.08162'
    'D14' = '92.91'
    'D15' = '5.298'
    'D17' = '14.82'
    'D18' = '0.000008525'
    'D21' = '4.983'
    'D22' = '10.62'
    'D23' = '6.390'
    'D24' = '2.284'
    'D25' = '146.56'
    'D26' = '1.738'
    'D27' = '18.03'
    'D28' = '114.05'
    'D29' = '4.710'
    'D30' = '4.608'
    'D31' = '0.09115'
    'D32' = '0.8112'
    'D33' = '0.05014'
    'D34' = '1.168'
    'D35' = '2.964'
    'D36' = '0.6072'
    'D37' = '2.653'
    'D38' = '3.203'
    'D39' = '0.01955'
    'D40' = '1.067'
    'D41' = '6.542'
    'D42' = '8.837'
    'D43' = '0.5161'
    'D44' = '114.94'
    'D45' = '0.1494'
    'D47' = '1.638'
    'D49' = '37.49'
    'D50' = '0.06043'
    'D51' = '62.25'
}

$plainUpdates = @{
    'D2' = '26.801.77'
    'E2' = '  -1.68%  '
    'D3' = '1.872.11'
    'E3' = '  -1.84%  '
    'E5' = '  -2.13%  '
    'E6' = '  +0.06%  '
    'E7' = '  +1.32%  '
    'E8' = '  -1.96%  '
    'E9' = '  -1.69%  '
    'E10' = '  -0.06%  '
    'E11' = '  -1.87%  '
    'E12' = '  -0.13%  '
    'D13' = '1.908.58'
    'E13' = '  +8.36%  '
    'E14' = '  -3.77%  '
    'E15' = '  -1.39%  '
    'E16' = '  +0.15%  '
    'E17' = '  +0.51%  '
    'E18' = '  -1.84%  '
    'E19' = '  +0.06%  '
    'D20' = '26.847.24'
    'E20' = '  -1.61%  '
    'E21' = '  -2.78%  '
    'E22' = '  -1.86%  '
    'E23' = '  -1.93%  '
    'E24' = '  -2.61%  '
    'E26' = '  -0.10%  '
    'E27' = '  -1.20%  '
    'E28' = '  -2.47%  '
    'E29' = '  -2.92%  '
    'E30' = '  -5.41%  '
    'E31' = '  -1.52%  '
    'E32' = '  -1.58%  '
    'E33' = '  -1.34%  '
    'E34' = '  -4.92%  '
    'E35' = '  -0.86%  '
    'E36' = '  +5.31%  '
    'E37' = '  -3.78%  '
    'E38' = '  -4.69%  '
    'E39' = '  -2.47%  '
    'E40' = '  -1.50%  '
    'E41' = '  -0.90%  '
    'E42' = '  -2.98%  '
    'E43' = '  +4.99%  '
    'E44' = '  -1.69%  '
    'E45' = '  -1.88%  '
    'E46' = '  +0.09%  '
    'E47' = '  -0.43%  '
    'E48' = '  -2.14%  '
    'E49' = '  -3.30%  '
    'E50' = '  -0.14%  '
    'E51' = '  -3.63%  '
}

foreach ($cellRef in $textUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $textUpdates[$cellRef]
}

foreach ($cellRef in $plainUpdates.Keys) {
    $ws.Range($cellRef).Value = $plainUpdates[$cellRef]
}
